$d = $word.ActiveDocument

$titleText  = "Play Alien Planets Free: Review of La Tuko's Cute and Entertaining Slot"
$oldBlurb   = "Looking for a fun space-themed slot game? Read our review of Alien Planets, play for free and discover its charming graphics and simple mechanics."
$newBlurb   = 'Create a cartoon-style feature image that showcases a happy Maya warrior with glasses playing "Alien Planets". The Maya warrior should be shown sitting in front of a computer screen, with the game''s logo and some of the adorable aliens from the slot displayed on the monitor. The image should have a vivid color palette, and the Maya warrior should have a joyful expression on his face conveying the excitement of the game. The background can be a depiction of an alien planet, with cute, funny-looking Martians scattered throughout, making it clear that the game has a space-themed backdrop. Overall, the image should convey the fun and easy gameplay of "Alien Planets", urging potential players to give it a spin.'
$metaRest   = ": Looking for a fun space-themed slot game? Read our review of Alien Planets, play for free and discover its charming graphics and simple mechanics."

# ---------------------------------------------------------------------------
# 1. Replace the wording of the final (italic) paragraph with the new
#    image-prompt text, keeping its existing italic run formatting intact.
#    (Done first, while the old blurb text is still unique in the document.)
# ---------------------------------------------------------------------------
$blurbRange = $d.Content
$blurbRange.Find.Execute($oldBlurb) | Out-Null
$blurbRange.Text = $newBlurb

Write-Output "step1 done"

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph that had been pasted again
#    near the bottom of the document.
# ---------------------------------------------------------------------------
$dupRange = $d.Content
$dupRange.Find.Execute($titleText) | Out-Null          # first hit: the real title
$dupRange.Find.Execute($titleText) | Out-Null          # second hit: the duplicate
$dupRange.Expand(4) | Out-Null                          # wdParagraph -> include mark
$dupRange.Delete()

Write-Output "step2 done"

# ---------------------------------------------------------------------------
# 3. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute($titleText) | Out-Null

$titlePara = $titleRange.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

# Locate the freshly-inserted (still empty) paragraph right after the title.
$metaPara = $d.Range($titleRange.End + 1, $titleRange.End + 1).Paragraphs(1)

# Put it back on the (implicit/default) Normal style instead of inheriting
# the title's Heading1 style.
$metaPara.Style = $d.Styles("Normal")

$metaRange = $metaPara.Range
$metaRange.Text = "Meta description" + $metaRest

$boldRange = $d.Range($metaRange.Start, $metaRange.Start + 16)
$boldRange.Bold = 1

# Mirror the leading empty run (<w:r/>) that precedes the text runs in every
# other body paragraph of this document.
$leadRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$leadRange.InsertBefore("")

Write-Output "step3 done"
